# Update gh-pages to output generated at 456a3b4
# Applies incremental "want to go" (F column) count bumps across the four
# sheets, and removes the expired "2023.01.12 上海·日漫咖啡体验" row from
# the 本地生活 (Local Life) sheet, shifting the remaining rows up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 展览 (Exhibitions) — column F (想去人数) count updates only
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1Updates = @{
    2  = 592
    4  = 6387
    5  = 722
    6  = 1087
    8  = 321
    10 = 17
    11 = 696
    12 = 1178
    13 = 79
    14 = 424
    16 = 19
    18 = 669
    19 = 383
    20 = 398
    22 = 1070
    23 = 145
    24 = 2214
    25 = 258
    26 = 96
    27 = 397
    29 = 3572
    30 = 47
}
foreach ($row in $ws1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $ws1Updates[$row]
}

# ---------------------------------------------------------------------
# 演出 (Performances) — column F (想去人数) count updates only
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2Updates = @{
    4  = 178
    8  = 707
    9  = 107
    13 = 102
    18 = 378
    19 = 314
    20 = 4088
    24 = 193
    25 = 232
    29 = 33
    32 = 1660
    34 = 2
}
foreach ($row in $ws2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $ws2Updates[$row]
}

# ---------------------------------------------------------------------
# 本地生活 (Local Life) — drop the expired first listing
# (2023.01.12 上海·日漫咖啡体验), shifting rows 3..11 up into 2..10,
# then remove the now-duplicate trailing row 11.
# Column A (serial index) is left untouched by this edit — it continues
# the existing 1..9 sequence.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3Rows = @{
    2  = @("2023.10.16", "上海·古影文化《1941·新和医院》大型沉浸式互动剧场", "金玉路2号 古影沉浸式互动游戏剧场", "2023.10.16 10:00-2024.10.15 21:00", 89, 996, "https://show.bilibili.com/platform/detail.html?id=77530", "//i0.hdslb.com/bfs/openplatform/202310/JqP3lHJt1698030195136.jpeg")
    3  = @("2023.10.25", "上海·方块大战（豫园店）", "丽水路88号2楼213 城隍庙第一购物中心", "2023.10.25 10:00-2024.10.20 21:00", 45, 49.9, "https://show.bilibili.com/platform/detail.html?id=79057", "//i1.hdslb.com/bfs/openplatform/202312/ASamaqBx1701419480253.jpeg")
    4  = @("2024.01.06", "上海·罗小黑 x HAPPY ZOO主题Cafe", "南京东路340号 百联zx创趣场", "2024.01.06 00:00-03.31 23:59", 1192, 10, "https://show.bilibili.com/platform/detail.html?id=80171", "//i2.hdslb.com/bfs/openplatform/202312/chPePM8d1703485388734.png")
    5  = @("2024.01.22", "上海·「新春特惠」世嘉都市乐园-JP国潮杂技嘉年华", "中山北路3300号环球港购物中心4楼 上海世嘉都市乐园", "2024.01.22 14:00-03.03 18:40", 2, "已停售", "https://show.bilibili.com/platform/detail.html?id=81210", "//i2.hdslb.com/bfs/openplatform/202401/sw2khwYM1706086166106.jpeg")
    6  = @("2024.01.27", "上海・明日方舟主题店·[SWEET ZONE甜蜜区域]", "南京东路830号第一百货商业中心B馆5楼(海底捞旁边) 第一百货商业中心", "2024.01.27 00:00-03.31 23:59", 1574, 30, "https://show.bilibili.com/platform/detail.html?id=81277", "//i0.hdslb.com/bfs/openplatform/202401/hp6D0Drt1705991831205.jpeg")
    7  = @("2024.02.01", "上海·次元波板糖×线条小狗MALTESE 主题快闪店", "西藏北路166静安大悦城北座6楼611号 次元波板糖", "2024.02.01 00:00-03.17 23:59", 430, 30, "https://show.bilibili.com/platform/detail.html?id=81345", "//i0.hdslb.com/bfs/openplatform/202401/Qbpful951706080847394.png")
    8  = @("2024.02.02", "上海·2024《永远的7日之都》x  萌果酱谷子咖啡", "南京东路340号百联ZX 萌果酱谷子咖啡（百联）", "2024.02.02 00:00-03.10 23:59", 129, 30, "https://show.bilibili.com/platform/detail.html?id=81357", "//i2.hdslb.com/bfs/openplatform/202401/5OYoWSGL1706087914805.jpeg")
    9  = @("2024.02.24", "上海·飘起来吧魔法泡泡-魔术表演（取消）", "曹杨路1888号 上海露边社·演艺空间", "2024.02.24 19:00-03.03 20:10", 5, "不可售", "https://show.bilibili.com/platform/detail.html?id=81524", "//i0.hdslb.com/bfs/openplatform/202401/tls18D0J1706599640356.png")
    10 = @("2024.03.01", "上海·「PLAVE with animate cafe」", "西藏北路198号大悦城北座8楼N809-1 animate cafe上海店", "2024.03.01 00:00-03.25 23:59", 769, 30, "https://show.bilibili.com/platform/detail.html?id=81873", "//i1.hdslb.com/bfs/openplatform/202402/7QENUAuN1708247451105.png")
}

foreach ($row in $ws3Rows.Keys) {
    $vals = $ws3Rows[$row]
    # Column B holds plain "yyyy.mm.dd" text in the source file; force a
    # text number-format first so Excel doesn't auto-coerce it to a date
    # serial the way it would for a bare Value assignment.
    $ws3.Cells.Item($row, 2).NumberFormat = "@"
    $ws3.Cells.Item($row, 2).Value = $vals[0]
    $ws3.Cells.Item($row, 3).Value = $vals[1]
    $ws3.Cells.Item($row, 4).Value = $vals[2]
    $ws3.Cells.Item($row, 5).Value = $vals[3]
    $ws3.Cells.Item($row, 6).Value = $vals[4]
    $ws3.Cells.Item($row, 7).Value = $vals[5]
    $ws3.Cells.Item($row, 8).Value = $vals[6]
    $ws3.Cells.Item($row, 9).Value = $vals[7]
}

# The old last row (11, "2024.03.01 PLAVE" before the shift) is now a
# duplicate of the new row 10 — remove it so the sheet has 9 data rows
# again (dimension A1:I11 -> A1:I10).
$ws3.Rows.Item(11).Delete()

# ---------------------------------------------------------------------
# 全部类型 (All Types) — column F (想去人数) count updates only.
# This sheet only mirrors on-sale listings, so the removed Local-Life row
# (which was already marked 已停售 there) never appeared here and no row
# shift is required.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4Updates = @{
    3  = 1192
    4  = 1574
    7  = 769
    8  = 592
    10 = 6387
    12 = 722
    13 = 1087
    14 = 707
    16 = 321
    18 = 696
    19 = 102
    21 = 1178
    22 = 79
    23 = 424
    25 = 378
    29 = 669
    30 = 383
    31 = 398
    33 = 193
    34 = 232
    35 = 1070
    36 = 145
    38 = 2214
    39 = 33
    40 = 1660
    41 = 258
    42 = 96
    43 = 397
    45 = 3572
    46 = 2
    48 = 47
}
foreach ($row in $ws4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $ws4Updates[$row]
}
